# Applies the "Substantial edits to Generic. Changed tools list in Master"
# commit to the CV document.
#
# Three content edits:
#  1. Reword the Selenium WebDriver / "studying" bullet.
#  2. Reword the "Tools and technologies ... include:" intro line.
#  3. Replace + italicise the long tools/technologies list, reordering and
#     adding/removing a few entries.
#  4. Drop two now-stale <w:lastRenderedPageBreak/> markers (one of which
#     disappears automatically because its whole paragraph is rewritten).

$d = $word.ActiveDocument

function Get-ParagraphRangeByText($doc, [string]$anchorText) {
    $probe = $doc.Content
    $found = $probe.Find.Execute($anchorText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        throw "Anchor text not found: $anchorText"
    }
    $para = $probe.Paragraphs(1)
    return $para.Range
}

$xmlHeader = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>'
$xmlFooter = '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

# ---------------------------------------------------------------------------
# 1. "I am experimenting with Selenium WebDriver as an automation tool for
#    regression checks, and am studying basic coding (Ruby)." paragraph.
#    Only the first three runs change; the " basic coding (Ruby)." tail is
#    left untouched (re-emitted verbatim, now merged into one run).
# ---------------------------------------------------------------------------
$seleniumPara = Get-ParagraphRangeByText $d "I am experimenting with Selenium WebDriver"

$seleniumBody = '<w:p w:rsidR="00011FB5" w:rsidRPr="009221A9" w:rsidRDefault="00011FB5" w:rsidP="00011FB5">' + `
  '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="18"/></w:numPr>' + `
  '<w:tabs><w:tab w:val="left" w:pos="0"/><w:tab w:val="left" w:pos="357"/><w:tab w:val="left" w:pos="720"/><w:tab w:val="left" w:pos="1440"/><w:tab w:val="left" w:pos="2160"/></w:tabs>' + `
  '<w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr></w:pPr>' + `
  '<w:r><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr>' + `
  '<w:t xml:space="preserve">I am experimenting with Selenium WebDriver as </w:t></w:r>' + `
  '<w:r><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr>' + `
  '<w:t>a tool for automated regression checks and am stud</w:t></w:r>' + `
  '<w:r><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr>' + `
  '<w:t>ying</w:t></w:r>' + `
  '<w:r><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr>' + `
  '<w:t xml:space="preserve"> basic coding (Ruby).</w:t></w:r>' + `
  '</w:p>'

$seleniumPara.InsertXML($xmlHeader + $seleniumBody + $xmlFooter)

# ---------------------------------------------------------------------------
# 2. "Tools and technologies I have explored ... include: " intro sentence.
#    Also drops the (now redundant) autoSpaceDN / textAlignment pPr settings.
# ---------------------------------------------------------------------------
$toolsIntroPara = Get-ParagraphRangeByText $d "Tools and technologies I have explored"

$toolsIntroBody = '<w:p w:rsidR="0021657F" w:rsidRPr="0021657F" w:rsidRDefault="00BA227F" w:rsidP="0021657F">' + `
  '<w:pPr><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/><w:overflowPunct/><w:autoSpaceDE/><w:adjustRightInd/>' + `
  '<w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/><w:noProof w:val="0"/><w:color w:val="222222"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr></w:pPr>' + `
  '<w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/><w:noProof w:val="0"/><w:color w:val="222222"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr>' + `
  '<w:t xml:space="preserve">Tools and technologies I have explored recently, either for work or for my own interest, include: </w:t></w:r>' + `
  '</w:p>'

$toolsIntroPara.InsertXML($xmlHeader + $toolsIntroBody + $xmlFooter)

# ---------------------------------------------------------------------------
# 3. The tools/technologies list paragraph: reworded, reordered, and
#    italicised. This also removes the stray <w:lastRenderedPageBreak/>
#    that used to sit mid-list.
# ---------------------------------------------------------------------------
$toolsListPara = Get-ParagraphRangeByText $d "Selenium WebDriver, Ruby, Firebug, Browser Developer Tools"

$toolsListBody = '<w:p w:rsidR="0021657F" w:rsidRPr="0021657F" w:rsidRDefault="0021657F" w:rsidP="00BA227F">' + `
  '<w:pPr><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/><w:overflowPunct/><w:autoSpaceDE/><w:adjustRightInd/><w:ind w:left="357"/>' + `
  '<w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/><w:i/><w:noProof w:val="0"/><w:color w:val="222222"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr></w:pPr>' + `
  '<w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/><w:i/><w:noProof w:val="0"/><w:color w:val="222222"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr>' + `
  '<w:t>Selenium WebDriver, Ruby, Browser Developer Tools, Rapid Reporter, VMWare, JIRA, Fiddler, RubyMine, RSpec, GitHub, TestLink, Sauce Labs, Applitools Eyes Express, Ubuntu, Firebug, Google Analytics, Google Search Console (Webmaster Tools), IntelliJ, Java.</w:t></w:r>' + `
  '</w:p>'

$toolsListPara.InsertXML($xmlHeader + $toolsListBody + $xmlFooter)

# ---------------------------------------------------------------------------
# 4. Drop the stale <w:lastRenderedPageBreak/> before the "Developing
#    brand-specific news ..." bullet (text itself is unchanged).
# ---------------------------------------------------------------------------
$devPara = Get-ParagraphRangeByText $d "Developing brand-specific news and information-based websites for the business."

$devBody = '<w:p w:rsidR="001307DD" w:rsidRPr="00CD3404" w:rsidRDefault="00BC272C" w:rsidP="00EB2615">' + `
  '<w:pPr><w:tabs><w:tab w:val="left" w:pos="0"/><w:tab w:val="left" w:pos="1440"/></w:tabs>' + `
  '<w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr></w:pPr>' + `
  '<w:r w:rsidRPr="00CD3404"><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr>' + `
  '<w:t>Developing brand-specific news and information-based websites for the business.</w:t></w:r>' + `
  '</w:p>'

$devPara.InsertXML($xmlHeader + $devBody + $xmlFooter)

Write-Host "Edits applied."
